$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$replacement = "MINDRAY MODELO DC " + [char]0x2013 + " N3 "

$find.Execute(
    "MEDISONIC MODELO H60 ",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    $replacement,
    2
)
